# The commit swaps the theme content between ppt/theme/theme1.xml (the
# theme used by the slide master / actual slides - previously the
# "Integral" colour palette) and ppt/theme/theme2.xml (the theme used by
# the notes master - previously the default "Office Theme" palette).
#
# fontScheme and fmtScheme are byte-identical between the two themes, so
# the only real content difference is the 12-colour clrScheme. We drive
# that through the SlideMaster's ColorScheme, which is the only theme
# colour surface PowerPoint's automation model exposes writes through.
#
# Office Theme's 12 standard theme colours (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink), expressed as COM RGB() integers (0x00BBGGRR):
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
